$d = $word.ActiveDocument

# --- Change 1: merge "6" + ".API COSTS" runs into a single run "6.API COSTS" ---
# The two runs share identical formatting, so a Find/Replace over the combined
# text collapses them into one run (matching the target XML).
$null = $d.Content.Find.Execute("6.API COSTS", $false, $false, $false, $false, $false, $true, 1, $false, "6.API COSTS", 2)

# --- Change 2: expand "All the pricings are estimated ." into several runs and
# relocate the "_GoBack" bookmark so it sits between "estimate" and "d, can" ---
$r = $d.Content
$found = $r.Find.Execute("All the pricings are estimated .", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Clear the matched text, leaving a collapsed insertion point.
    $r.Text = ""

    $p1 = $d.Range($r.Start, $r.Start)
    $p1.InsertAfter("*")

    $p2 = $d.Range($p1.End, $p1.End)
    $p2.InsertAfter("all prices are estimate")

    $p3 = $d.Range($p2.End, $p2.End)
    $p3.InsertAfter("d, can")

    $p4 = $d.Range($p3.End, $p3.End)
    $p4.InsertAfter(" vary while actual implementation")

    $p5 = $d.Range($p4.End, $p4.End)
    $p5.InsertAfter(".")

    # Move the existing "_GoBack" bookmark to sit right after "estimate"
    # (re-adding with the same name relocates it rather than duplicating it).
    $bmRange = $d.Range($p2.End, $p2.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
